# A new weekly price-report row (Jengibre / Vega Modelo de Temuco) is
# inserted at row 325, pushing all the existing rows 325-359 down to
# 326-360 (dimension grows from A1:R359 to A1:R360).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 325, shifting rows down.
$ws.Rows.Item(325).Insert()

# Populate the newly inserted row with the latest weekly report entry.
$ws.Cells.Item(325, 1).Value2 = 10
$ws.Cells.Item(325, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(325, 3).Value2 = "La Araucanía"
$ws.Cells.Item(325, 4).Value2 = 45194
$ws.Cells.Item(325, 5).Value2 = 9
$ws.Cells.Item(325, 6).Value2 = 100114007
$ws.Cells.Item(325, 7).Value2 = "Jengibre"
$ws.Cells.Item(325, 8).Value2 = "Sin especificar"
$ws.Cells.Item(325, 9).Value2 = "Primera"
$ws.Cells.Item(325, 10).Value2 = 160
$ws.Cells.Item(325, 11).Value2 = 22000
$ws.Cells.Item(325, 12).Value2 = 24000
$ws.Cells.Item(325, 13).Value2 = 23000
$ws.Cells.Item(325, 14).Value2 = "$/caja 13 kilos"
$ws.Cells.Item(325, 15).Value2 = "Perú"
$ws.Cells.Item(325, 16).Value2 = 1769
$ws.Cells.Item(325, 17).Value2 = 13
$ws.Cells.Item(325, 18).Value2 = "Hortaliza"
